# Edit script: update holdings weights/percent-change values and the
# "as of" date in the confidential disclaimer string, per the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet is protected; temporarily unprotect so the cell values can
# be updated, then restore protection afterwards.
$ws.Unprotect()

# Update the "Weight" (D) and "Percent Change" (E) columns for rows 2-77.
$ws.Range("D2").Value2 = 0.07636341518540066
$ws.Range("E2").Value2 = -0.01999347151950381
$ws.Range("D3").Value2 = 0.04692480139922205
$ws.Range("E3").Value2 = -0.0160733067729083
$ws.Range("D4").Value2 = 0.03671715525197795
$ws.Range("E4").Value2 = -0.008923310042932897
$ws.Range("D5").Value2 = 0.03320655693598281
$ws.Range("E5").Value2 = -0.000384338217631619
$ws.Range("D6").Value2 = 0.03129476699276117
$ws.Range("E6").Value2 = 0.007761273919443257
$ws.Range("D7").Value2 = 0.03053036010845385
$ws.Range("E7").Value2 = -0.004310914942708854
$ws.Range("D8").Value2 = 0.03037735908269638
$ws.Range("E8").Value2 = 0.00972871842843781
$ws.Range("D9").Value2 = 0.02847209997504797
$ws.Range("E9").Value2 = 0.005315935725504373
$ws.Range("D10").Value2 = 0.0265339374210001
$ws.Range("E10").Value2 = -0.008560112431327482
$ws.Range("D11").Value2 = 0.0277645263163397
$ws.Range("E11").Value2 = -0.0003016591251884249
$ws.Range("D12").Value2 = 0.02446311415351443
$ws.Range("E12").Value2 = -0.02124728448047464
$ws.Range("D13").Value2 = 0.02317898237721572
$ws.Range("E13").Value2 = 0
$ws.Range("D14").Value2 = 0.02046654190494069
$ws.Range("E14").Value2 = -0.02706731752944946
$ws.Range("D15").Value2 = 0.0186062009869982
$ws.Range("E15").Value2 = -0.01714824659178604
$ws.Range("D16").Value2 = 0.02069342612397693
$ws.Range("E16").Value2 = -0.007198560287942413
$ws.Range("D17").Value2 = 0.01869673531830337
$ws.Range("E17").Value2 = -0.00619680453934579
$ws.Range("D18").Value2 = 0.01807675462220698
$ws.Range("E18").Value2 = -0.00466083464792777
$ws.Range("D19").Value2 = 0.01514106925148231
$ws.Range("E19").Value2 = 0.02028250633828343
$ws.Range("D20").Value2 = 0.01411198903848604
$ws.Range("E20").Value2 = -0.01697825257535301
$ws.Range("D21").Value2 = 0.01593788104512244
$ws.Range("E21").Value2 = -0.02921240064687058
$ws.Range("D22").Value2 = 0.01416922109373843
$ws.Range("E22").Value2 = 0.01232865144819439
$ws.Range("D23").Value2 = 0.01360587420841443
$ws.Range("E23").Value2 = -0.02842632898034558
$ws.Range("D24").Value2 = 0.01525523423974709
$ws.Range("E24").Value2 = -0.0003333333333332966
$ws.Range("D25").Value2 = 0.01469687272508968
$ws.Range("E25").Value2 = 0.005698778833107321
$ws.Range("D26").Value2 = 0.0119210183379357
$ws.Range("E26").Value2 = -0.02810304449648726
$ws.Range("D27").Value2 = 0.01203767601153376
$ws.Range("E27").Value2 = -0.002981860349540177
$ws.Range("D28").Value2 = 0.01234238186667538
$ws.Range("E28").Value2 = -0.02268431001890359
$ws.Range("D29").Value2 = 0.01231466320576917
$ws.Range("E29").Value2 = -0.00391877449234046
$ws.Range("D30").Value2 = 0.01085634257838002
$ws.Range("E30").Value2 = 0.01476828125861029
$ws.Range("D31").Value2 = 0.01192600370860228
$ws.Range("E31").Value2 = 0.00953097567093053
$ws.Range("D32").Value2 = 0.0130109200730643
$ws.Range("E32").Value2 = 0.009552382922960234
$ws.Range("D33").Value2 = 0.01126718697501321
$ws.Range("E33").Value2 = -0.0146899404880424
$ws.Range("D34").Value2 = 0.01163241523004715
$ws.Range("E34").Value2 = 0.001757160428747184
$ws.Range("D35").Value2 = 0.009903339121755816
$ws.Range("E35").Value2 = -0.04816056542225444
$ws.Range("D36").Value2 = 0.01101652253789736
$ws.Range("E36").Value2 = 0.002529675034053325
$ws.Range("D37").Value2 = 0.01080479384568754
$ws.Range("E37").Value2 = 0.01674895030683365
$ws.Range("D38").Value2 = 0.009940579840635198
$ws.Range("E38").Value2 = 0.01143458963364163
$ws.Range("D39").Value2 = 0.009254044446139913
$ws.Range("E39").Value2 = -0.01463172865577733
$ws.Range("D40").Value2 = 0.009610996985867329
$ws.Range("E40").Value2 = -0.0204996265250228
$ws.Range("D41").Value2 = 0.009249258490299993
$ws.Range("E41").Value2 = -0.04029580440688207
$ws.Range("D42").Value2 = 0.009570316361228003
$ws.Range("E42").Value2 = -0.02887981330221712
$ws.Range("D43").Value2 = 0.009886089739249435
$ws.Range("E43").Value2 = -0.0003933394519470834
$ws.Range("D44").Value2 = 0.009390244772750989
$ws.Range("E44").Value2 = -0.006540805708339481
$ws.Range("D45").Value2 = 0.009150448443688298
$ws.Range("E45").Value2 = 0.02677258017063844
$ws.Range("D46").Value2 = 0.009519266165602183
$ws.Range("E46").Value2 = 0.00707013574660631
$ws.Range("D47").Value2 = 0.008696081761135831
$ws.Range("E47").Value2 = -0.001651073197578512
$ws.Range("D48").Value2 = 0.007325703072305216
$ws.Range("E48").Value2 = -0.01910932055749137
$ws.Range("D49").Value2 = 0.008271826717409531
$ws.Range("E49").Value2 = -0.002772387025228751
$ws.Range("D50").Value2 = 0.007945434499868275
$ws.Range("E50").Value2 = -0.01204705882352941
$ws.Range("D51").Value2 = 0.00761640003587373
$ws.Range("E51").Value2 = 0.0136540664375715
$ws.Range("D52").Value2 = 0.007527411169475206
$ws.Range("E52").Value2 = 0.01821974965229467
$ws.Range("D53").Value2 = 0.006825969516686834
$ws.Range("E53").Value2 = -0.0006573181419807783
$ws.Range("D54").Value2 = 0.007144734117108217
$ws.Range("E54").Value2 = 0.008624419107693626
$ws.Range("D55").Value2 = 0.006648941497001769
$ws.Range("E55").Value2 = -0.001715165550028264
$ws.Range("D56").Value2 = 0.006375690838080966
$ws.Range("E56").Value2 = 0.0003127736769672484
$ws.Range("D57").Value2 = 0.006836737917326656
$ws.Range("E57").Value2 = -0.005250262513125659
$ws.Range("D58").Value2 = 0.006363227411414506
$ws.Range("E58").Value2 = -0.01170497814130578
$ws.Range("D59").Value2 = 0.005498165893348848
$ws.Range("E59").Value2 = 0.01865150608418076
$ws.Range("D60").Value2 = 0.006239689926296554
$ws.Range("E60").Value2 = 0.001438159156280205
$ws.Range("D61").Value2 = 0.00531440513057856
$ws.Range("E61").Value2 = 0.02026266416510336
$ws.Range("D62").Value2 = 0.005842056761929813
$ws.Range("E62").Value2 = 0.002594210813763098
$ws.Range("D63").Value2 = 0.005395367550203885
$ws.Range("E63").Value2 = 0.002956830277942135
$ws.Range("D64").Value2 = 0.004994942578263857
$ws.Range("E64").Value2 = -0.00511018843819866
$ws.Range("D65").Value2 = 0.004791140625413902
$ws.Range("E65").Value2 = 0.006284858070423782
$ws.Range("D66").Value2 = 0.004328597934968237
$ws.Range("E66").Value2 = -0.001174763319742866
$ws.Range("D67").Value2 = 0.004391712727607191
$ws.Range("E67").Value2 = 0.002179539572265243
$ws.Range("D68").Value2 = 0.003645851422179557
$ws.Range("E68").Value2 = -0.03897116134060796
$ws.Range("D69").Value2 = 0.003990440242653844
$ws.Range("E69").Value2 = 0.005547018477568333
$ws.Range("D70").Value2 = 0.003697250593752038
$ws.Range("E70").Value2 = 0.01275585879560959
$ws.Range("D71").Value2 = 0.00321556407994669
$ws.Range("E71").Value2 = -0.02544186046511621
$ws.Range("D72").Value2 = 0.002606501345610121
$ws.Range("E72").Value2 = -0.03272574259319472
$ws.Range("D73").Value2 = 0.002667621989982441
$ws.Range("E73").Value2 = -0.02668709936646185
$ws.Range("D74").Value2 = 0.002294267580761964
$ws.Range("E74").Value2 = -0.01888309430682311
$ws.Range("D75").Value2 = 0.001850968421089313
$ws.Range("E75").Value2 = 0.01599870717517793
$ws.Range("D76").Value2 = 0.001767912145784024
$ws.Range("E76").Value2 = -0.01894986182392411
$ws.Range("D77").Value2 = 1
$ws.Range("E77").Value2 = -0.006100485913858922

# Update the "as of" date in the confidential disclaimer text (A80) from
# 2021-03-23 to 2021-03-24.
$ws.Range("A80").Value2 = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-24 for illustrative purposes only and are subject to change."

# Re-fit the row so updating the wrapped text doesn't leave a stray custom
# row-height behind (keeps row 80 identical in shape to the original file).
$ws.Rows.Item(80).AutoFit()

# Restore sheet protection to match the original workbook's protected state.
$ws.Protect()
